# Generate Report for Handback
# The "6135c9e7-..." file has now been handed back (in sync with en-US)
# for both locales, so update the Overview sheet and the per-locale
# report sheets accordingly, and clear the stale "out of date" error.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the 6135c9e7-... file ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = $status   # zh-cn column
$ws.Range("F3").Value = $status   # de-de column

# --- zh-cn report sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = $status                      # Status
$ws.Range("K3").Value = "2016-09-06 21:06:47"         # Latest Handback DateTime
$ws.Range("P3").Value = ""                            # Error Detail cleared
$ws.Columns.Item(16).ColumnWidth = 13.7470528738839

# --- de-de report sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = $status                      # Status
$ws.Range("K3").Value = "2016-09-06 21:06:55"         # Latest Handback DateTime
$ws.Range("P3").Value = ""                            # Error Detail cleared
$ws.Columns.Item(16).ColumnWidth = 13.7470528738839
